$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: remove the block of empty placeholder rows 18-26 (without shifting
# the rows below, since they keep their original row numbers 27-30), then
# append 9 more identical empty placeholder rows after row 30 (31-39). ---
$ws1.Rows("18:26").Clear()

$srcRow = $ws1.Range("A27:C27")
for ($r = 31; $r -le 39; $r++) {
    $dstRow = $ws1.Range("A" + $r + ":C" + $r)
    $srcRow.Copy($dstRow)
}

# --- Sheet2: remove the duplicated data table (rows 1-10) that mirrored
# Sheet1's contents, keeping the single trailing placeholder row, whose
# row number shifts from 34 to 33. ---
$ws2.Rows("1:10").Delete()
$ws2.Rows("1:9").Insert()

# Update the selections to match the post-edit state.
$ws2.Range("A1:XFD9").Select()
$ws1.Select()
$ws1.Range("A16").Select()
